# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Espárragos" at the top of the
# existing data block (row 20), pushing the previous rows 20-52 down to
# rows 22-54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 20 (each Insert() shifts
# everything at/below down by one row).
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# New row 20
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 45203
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = "Espárragos"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 1300
$ws.Range("L20").Value = 1300
$ws.Range("M20").Value = 1300
$ws.Range("N20").Value = "$/kilo"
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 1300
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"

# New row 21
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 45203
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 300000000
$ws.Range("G21").Value = "Espárragos"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1500
$ws.Range("N21").Value = "$/kilo"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1500
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
